$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G ("K") rows 2-72. The analysis script was
# regenerated to compute K from the updated source ("K instead of
# Strike#") and the std/mean pipeline + s_vals were recalculated and
# written back into this save_data sheet.
$kValues = @(
    0,0,0,0,2,1,1,1,0,0,1,2,1,2,2,1,2,0,1,1,4,1,1,0,0,1,1,1,1,2,0,0,2,2,2,1,0,0,1,1,2,2,0,1,3,
    1,0,1,0,1,2,1,4,3,1,0,1,1,0,0,1,0,1,1,2,0,1,2,1,1,2
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
